# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" (Exhibition) sheet and the "全部类型" (All Types) sheet.
#
# Mapping of changes (row -> old -> new), identical on both sheets except
# that the last updated row differs (row 12 on 展览, row 14 on 全部类型):
#   F2:  268  -> 269
#   F5:  6666 -> 6675
#   F6:  5447 -> 5468
#   F9:  7    -> 8
#   F11: 237  -> 239
#   F(last): 100 -> 116

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 269
$wsExhibition.Range("F5").Value = 6675
$wsExhibition.Range("F6").Value = 5468
$wsExhibition.Range("F9").Value = 8
$wsExhibition.Range("F11").Value = 239
$wsExhibition.Range("F12").Value = 116

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 269
$wsAllTypes.Range("F5").Value = 6675
$wsAllTypes.Range("F6").Value = 5468
$wsAllTypes.Range("F9").Value = 8
$wsAllTypes.Range("F11").Value = 239
$wsAllTypes.Range("F14").Value = 116
